$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "LIKE" operator description (row 9, column B) to append
# the new sentence about bracket patterns, matching the commit's note update.
$ws.Range("B9").Value = "Search for a pattern. % and _ are wildcards. ``%`` is like ``+`` In regex(0/1/1+).  ``_`` is like ``.``. Also has ``[abc]`` and ``[^abc]`` like regex."

# Update the active selection to B10, as recorded in the saved view state.
$ws.Range("B10").Select()
